# Update cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '66.058.37'
$ws.Range('E2').Value = '  +2.59%  '

# Row 3
$ws.Range('D3').Value = '3.202.97'
$ws.Range('E3').Value = '  +5.58%  '

# Row 4
$ws.Range('E4').Value = '  +0.24%  '

# Row 5
$ws.Range('D5').Value = '''573.27'
$ws.Range('E5').Value = '  +4.16%  '

# Row 6
$ws.Range('D6').Value = '''150.39'
$ws.Range('E6').Value = '  +9.02%  '

# Row 7
$ws.Range('D7').Value = '''1.00'
$ws.Range('E7').Value = '  -0.09%  '

# Row 8
$ws.Range('D8').Value = '3.192.12'
$ws.Range('E8').Value = '  +5.49%  '

# Row 9
$ws.Range('D9').Value = '''0.510'
$ws.Range('E9').Value = '  +5.12%  '

# Row 10
$ws.Range('D10').Value = '''6.99'
$ws.Range('E10').Value = '  +10.08%  '

# Row 11
$ws.Range('D11').Value = '''0.161'
$ws.Range('E11').Value = '  +5.62%  '

# Row 12
$ws.Range('D12').Value = '''0.483'
$ws.Range('E12').Value = '  +6.00%  '

# Row 13
$ws.Range('D13').Value = '''38.14'
$ws.Range('E13').Value = '  +7.39%  '

# Row 14
$ws.Range('D14').Value = '''0.0000231'
$ws.Range('E14').Value = '  +5.98%  '

# Row 15
$ws.Range('D15').Value = '3.715.69'
$ws.Range('E15').Value = '  +5.69%  '

# Row 16
$ws.Range('D16').Value = '66.201.75'
$ws.Range('E16').Value = '  +2.86%  '

# Row 17
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.220.37'
$ws.Range('E17').Value = '  +5.97%  '

# Row 18
$ws.Range('B18').Value = 'BitcoinCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D18').Value = '''533.56'
$ws.Range('E18').Value = '  +11.01%  '

# Row 19
$ws.Range('E19').Value = '  +2.79%  '

# Row 20
$ws.Range('D20').Value = '''7.08'
$ws.Range('E20').Value = '  +7.80%  '

# Row 21
$ws.Range('D21').Value = '''14.49'
$ws.Range('E21').Value = '  +7.01%  '

# Row 22
$ws.Range('D22').Value = '''0.738'
$ws.Range('E22').Value = '  +8.50%  '

# Row 23
$ws.Range('D23').Value = '''7.68'
$ws.Range('E23').Value = '  +8.76%  '

# Row 24
$ws.Range('D24').Value = '''13.46'
$ws.Range('E24').Value = '  +8.84%  '

# Row 25
$ws.Range('D25').Value = '''80.72'
$ws.Range('E25').Value = '  +3.24%  '

# Row 26
$ws.Range('D26').Value = '''0.997'
$ws.Range('E26').Value = '  -0.19%  '

# Row 27
$ws.Range('D27').Value = '''9.35'
$ws.Range('E27').Value = '  +21.96%  '

# Row 28
$ws.Range('D28').Value = '''2.93'
$ws.Range('E28').Value = '  +8.62%  '

# Row 29
$ws.Range('D29').Value = '''2.24'
$ws.Range('E29').Value = '  +8.44%  '

# Row 30
$ws.Range('D30').Value = '''27.31'
$ws.Range('E30').Value = '  +6.45%  '

# Row 31
$ws.Range('B31').Value = 'FirstDigitalUSD'
$ws.Range('C31').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D31').Value = '''1.00'
$ws.Range('E31').Value = '  +0.15%  '

# Row 32
$ws.Range('B32').Value = 'Stacks'
$ws.Range('C32').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D32').Value = '''2.73'
$ws.Range('E32').Value = '  +6.23%  '

# Row 33
$ws.Range('D33').Value = '''1.17'
$ws.Range('E33').Value = '  +5.33%  '

# Row 34
$ws.Range('D34').Value = '''559.29'
$ws.Range('E34').Value = '  +0.91%  '

# Row 35
$ws.Range('D35').Value = '''6.31'
$ws.Range('E35').Value = '  +8.49%  '

# Row 36
$ws.Range('D36').Value = '''5.58'
$ws.Range('E36').Value = '  +5.28%  '

# Row 37
$ws.Range('D37').Value = '''54.74'
$ws.Range('E37').Value = '  +4.40%  '

# Row 38
$ws.Range('D38').Value = '''0.0446'
$ws.Range('E38').Value = '  +9.20%  '

# Row 39
$ws.Range('D39').Value = '''0.0852'
$ws.Range('E39').Value = '  +8.10%  '

# Row 40
$ws.Range('D40').Value = '''0.128'
$ws.Range('E40').Value = '  +7.40%  '

# Row 41
$ws.Range('D41').Value = '3.208.20'
$ws.Range('E41').Value = '  +10.15%  '

# Row 42
$ws.Range('D42').Value = '''2.90'
$ws.Range('E42').Value = '  +6.43%  '

# Row 43
$ws.Range('D43').Value = '''8.54'
$ws.Range('E43').Value = '  +4.92%  '

# Row 44
$ws.Range('D44').Value = '''0.282'
$ws.Range('E44').Value = '  +17.85%  '

# Row 45
$ws.Range('D45').Value = '''2.34'
$ws.Range('E45').Value = '  +13.54%  '

# Row 46
$ws.Range('D46').Value = '''26.34'
$ws.Range('E46').Value = '  +7.43%  '

# Row 47
$ws.Range('D47').Value = '''0.999'
$ws.Range('E47').Value = '  +0.02%  '

# Row 48
$ws.Range('D48').Value = '0.0₃0549'
$ws.Range('E48').Value = '  +5.12%  '

# Row 49
$ws.Range('D49').Value = '''124.43'
$ws.Range('E49').Value = '  +5.28%  '

# Row 50
$ws.Range('E50').Value = '  +4.29%  '

# Row 51
$ws.Range('D51').Value = '''2.19'
$ws.Range('E51').Value = '  +8.81%  '

